$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contributor row: name / email / repo link
$ws.Range("A3").Value = "مريم محمد قيس عثمان"
$ws.Range("B3").Value = "mariamqaies2020@gmail.com"
$ws.Range("C3").Value = "https://github.com/mariamqaies/Security-Task.git"

# Hyperlink the new e-mail address (mirrors the existing B2 mailto link)
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:mariamqaies2020@gmail.com")
$ws.Range("B3").Style = "Hyperlink"

# Leave the selection where the user ended up after entering the data
$ws.Range("B6").Select() | Out-Null
